$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10 (Q0..Q8), columns B..G
$data = @(
    @(2,  0.1374531038001369, 0.9463274094814685, 4.340354840748894, 2.083351828364306, 2.099094096015294, 52),
    @(3,  0.1201387303292092, 0.9941924500794057, 4.368983621713564, 2.090211382064877, 2.10752018100236,  51),
    @(4,  0.1570212086422995, 0.9129929970791744, 4.115459563428113, 2.028659548428004, 2.043107852613974, 50),
    @(5,  0.1668826930103534, 1.018518659023897,  4.584217853051491, 2.141078665778418, 2.156685501776154, 49),
    @(6,  0.1498087523549735, 0.9654930794235904, 4.38482909777947,  2.093998351904669, 2.110735219814999, 48),
    @(7,  0.1482159292065307, 1.037301259063642,  5.28735183827951,  2.299424240604484, 2.327192174260785, 36),
    @(8,  0.1958133101793198, 1.063338432913966,  5.471053399423887, 2.339028302399073, 2.364845904346298, 35),
    @(9,  0.1080685240234339, 1.516053278149146,  9.891401411246379, 3.145059842236134, 3.234328795688273, 18),
    @(10, -0.8008619424696254, 1.150395842953777, 6.142350263041148, 2.478376537784594, 2.459891756431453, 11)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
